$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 ---
$ws.Range("A8").Value = ""
$ws.Range("C8").Value = 40
$ws.Range("D8").Value = "'1.0"
$ws.Range("E8").Value = 'Rewiring of light point/ fan point/ exhaust fan point/ call bell point with 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade (IS:694) in recessed  ISI marked MMS ( IS:9537 P - III ) virgin material PVC conduit & it''s ISI marked (IS:3419-1988) accessories, round tiles, 1.2 mm thick MS box with earth terminal, 6 A switch, 3 pin ceiling rose/holder / 3 way connector , 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/  brass  screws, cup washers, making connections, testing etc. as required. For specification of copper  Conductor,  Phenolic Laminated sheet''s & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = "'0.00"

# --- Row 9 ---
$ws.Range("A9").Value = 'P. point'
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = "'3"
$ws.Range("E9").Value = 'Medium point (up to 6 mtr.)'
$ws.Range("F9").Value = 472
$ws.Range("G9").Value = "'472.00"

# --- Row 10 ---
$ws.Range("A10").Value = 'P. point'
$ws.Range("C10").Value = 87
$ws.Range("D10").Value = "'4"
$ws.Range("E10").Value = 'Long point  (up to 10 mtr.)'
$ws.Range("F10").Value = 662
$ws.Range("G10").Value = "'57594.00"

# --- Row 11 ---
$ws.Range("A11").Value = ""
$ws.Range("C11").Value = 51
$ws.Range("D11").Value = "'17.0"
$ws.Range("E11").Value = 'Providing & Fixing of 240/415 V AC MCB with positive isolation of 10 kA breaking capacity (B/ C/D tripping characteristic as per type of load and  site requirement) 4 KV impulse withstand voltage, ISI marked IS 8828(1996) / conforming to IEC 60898-1 2002, IEC 60947-2, low watt losses, trip free mechanisum , energy limiting of  class 3 as per IEC,  minimum phase termination capacity of 35sq.mm. , conductor line load reversibility , IP 20 contact protection and fitted in  existing distribution board/sheets, minimum electrical operation 20,000 upto 20 A rating and 10,000 upto 63 A, 5000 for 80 A & above rating  including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = "'0.00"

# --- Row 12 ---
$ws.Range("C12").Value = 51
$ws.Range("D12").Value = "'34"
$ws.Range("E12").Value = 'Metal door (single phase) IK-09 and IP-43 with Metal end box'

# --- Row 14 (Grand Total) ---
$ws.Range("G14").Value = "'58066.00"
$ws.Range("H14").Value = "'58066.00"

# --- Row 16 (Net Payable Amount) ---
$ws.Range("G16").Value = "'58066.00"
$ws.Range("H16").Value = "'58066.00"

